$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date / Contact, insert a Jurisdiction row ---
$ws = $wb.Worksheets.Item(1)

# Update the "Date" property value (row 8, column B)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update the "Contact" property value (row 10, column B)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row after "Contact" (row 10) for the "Jurisdiction" property,
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (border/alignment style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# --- Rename the code-system include sheet ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
